$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring formatting for the two new rows in line with the existing table
# (column A uses the "wrap text" style, column B the plain one) without
# minting any new style records.
$ws.Range("A1").Copy()
$ws.Range("A5:A6").PasteSpecial(-4122)

$ws.Range("B1").Copy()
$ws.Range("B5:B6").PasteSpecial(-4122)

$ws.Range("A5").Value = "valid"
$ws.Range("B5").Value = "invalid"
$ws.Range("A6").Value = "invalid"
$ws.Range("B6").Value = "invalid"

$ws.Range("A7").Select()
